# Generate Report for Handback
# Adds a new handback-status row for "e9857952-1f67-43b3-9185-4783f0f601d2.md"
# (alongside the existing "864724d1-6658-432c-854f-423524d4c2ce.md", which is
# the renamed/rerun version of the old "e71a5e0b-2db9-4931-9633-2e9f660ff036.md")
# across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$guid1 = "864724d1-6658-432c-854f-423524d4c2ce"
$guid2 = "e9857952-1f67-43b3-9185-4783f0f601d2"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Drop hyperlinks that need to be replaced before adding any new ones -- this
# engine's Hyperlinks.Delete() clears *all* links on the sheet, not just the
# target range, so every delete must happen before any add.
$wsOverview.Range("B2").Hyperlinks.Delete()

# Row 2 already describes guid1's file; refresh its hyperlink + generate date.
$wsOverview.Range("G2").Value = "2016-09-01 11:10:59"
$wsOverview.Range("G2").NumberFormat = $dateFmt

# New row 3 for guid2's file.
$wsOverview.Range("A3").Value = "$guid2.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-09-01 11:10:59"
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e582daa4816dd8057226d84a1c0d765367c22b1c/e2e/$guid1.md", $null, $null, "e2e\$guid1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e582daa4816dd8057226d84a1c0d765367c22b1c/e2e/$guid2.md", $null, $null, "e2e\$guid2.md")

$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Range("I2").Hyperlinks.Delete()

# Row 2: refresh guid1 file's handoff/handback info.
$wsZhCn.Range("G2").Value = "$guid1.f1abcc83ac2cb36839928c5bb2d95fc77c794265.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-01 11:10:54"
$wsZhCn.Range("H2").NumberFormat = $dateFmt
$wsZhCn.Range("J2").Value = "$guid1.f1abcc83ac2cb36839928c5bb2d95fc77c794265.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-01 11:11:21"
$wsZhCn.Range("K2").NumberFormat = $dateFmt

# Row 3: new guid2 file.
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'True"
$wsZhCn.Range("G3").Value = "$guid2.8edd37c7a623868290484a09e4f9839b446ead25.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-01 11:10:54"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("J3").Value = "$guid2.8edd37c7a623868290484a09e4f9839b446ead25.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-01 11:11:21"
$wsZhCn.Range("K3").NumberFormat = $dateFmt
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("O3").Value = "'False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e582daa4816dd8057226d84a1c0d765367c22b1c/e2e/$guid1.md", $null, $null, "$guid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3d669a4cf15c257a8172365f74ba25859c1d2139/e2e/$guid1.md", $null, $null, "$guid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e582daa4816dd8057226d84a1c0d765367c22b1c/e2e/$guid2.md", $null, $null, "$guid2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3d669a4cf15c257a8172365f74ba25859c1d2139/e2e/$guid2.md", $null, $null, "$guid2.md")

$loZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Range("I2").Hyperlinks.Delete()

# Row 2: refresh guid1 file's handoff/handback info.
$wsDeDe.Range("G2").Value = "$guid1.f1abcc83ac2cb36839928c5bb2d95fc77c794265.de-de.xlf"
$wsDeDe.Range("J2").Value = "$guid1.f1abcc83ac2cb36839928c5bb2d95fc77c794265.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-01 11:11:29"
$wsDeDe.Range("K2").NumberFormat = $dateFmt

# Row 3: new guid2 file.
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'True"
$wsDeDe.Range("G3").Value = "$guid2.8edd37c7a623868290484a09e4f9839b446ead25.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-01 11:10:59"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("J3").Value = "$guid2.8edd37c7a623868290484a09e4f9839b446ead25.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-01 11:11:29"
$wsDeDe.Range("K3").NumberFormat = $dateFmt
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("O3").Value = "'False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e582daa4816dd8057226d84a1c0d765367c22b1c/e2e/$guid1.md", $null, $null, "$guid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/702e0218d3bde65a58cc5867aa8b9c6a783a97e1/e2e/$guid1.md", $null, $null, "$guid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e582daa4816dd8057226d84a1c0d765367c22b1c/e2e/$guid2.md", $null, $null, "$guid2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/702e0218d3bde65a58cc5867aa8b9c6a783a97e1/e2e/$guid2.md", $null, $null, "$guid2.md")

$loDeDe = $wsDeDe.ListObjects.Item("de-de")
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
